# "testing CM with RS"
# Apply edits to the "Coupling Parameters" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year (row 3) changes from 2089 to 2055
$ws.Range("B3").Value = 2055

# scenarioWeatheryearsExcel (row 29) switches weather-year file
$ws.Range("B29").Value = "40weatherYears2050TNO.xlsx"

# reliability_option_SP (row 44): was a single "NOTSET" cell in B44.
# Now B44 holds a numeric value (150) and the "NOTSET" label moves to C44.
$ws.Range("B44").Value = 150
$ws.Range("C44").HorizontalAlignment = -4108
$ws.Range("C44").Value = "NOTSET"

# Update the view: scroll back to the top and select B4 instead of C30
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
